# Refresh Leve profit-tracker market data (currentAveragePrice* / LevePrice* /
# LeveProfit* columns H:N) on each job sheet to the latest Universalis pull.
# Mirrors the scheduled-runner job that keeps Hades_Profits.xlsx up to date.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 64: Forged from the Void (Void Glue)
$ws.Range("H64").Value = 4235.0625
$ws.Range("I64").Value = 3750
$ws.Range("J64").Value = 4526.1
$ws.Range("K64").Value = 3750
$ws.Range("L64").Value = 4526.1
$ws.Range("M64").Value = -3502
$ws.Range("N64").Value = -5022.1

# Row 67: Dodging the Draft (L) (Void Glue)
$ws.Range("H67").Value = 4235.0625
$ws.Range("I67").Value = 3750
$ws.Range("J67").Value = 4526.1
$ws.Range("K67").Value = 3750
$ws.Range("L67").Value = 4526.1
$ws.Range("M67").Value = -2892
$ws.Range("N67").Value = -6242.1

# Row 74: Adhesive of Antipathy (Wing Glue)
$ws.Range("H74").Value = 3903.4119
$ws.Range("I74").Value = 3817.5557
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3817.5557
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2881.5557
$ws.Range("N74").Value = -5872

# Row 76: Warding Off Temptation (Enchanted Hardsilver Ink)
$ws.Range("H76").Value = 3748.7368
$ws.Range("I76").Value = 3532.2222
$ws.Range("J76").Value = 3943.6
$ws.Range("K76").Value = 3532.2222
$ws.Range("L76").Value = 3943.6
$ws.Range("M76").Value = -3217.2222
$ws.Range("N76").Value = -4573.6

# Row 77: It's Gonna Grow Back (L) (Wing Glue)
$ws.Range("H77").Value = 3903.4119
$ws.Range("I77").Value = 3817.5557
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19087.7785
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14407.7785
$ws.Range("N77").Value = -29360

# Row 79: The Garden of Arcane Delights (L) (Enchanted Hardsilver Ink)
$ws.Range("H79").Value = 3748.7368
$ws.Range("I79").Value = 3532.2222
$ws.Range("J79").Value = 3943.6
$ws.Range("K79").Value = 3532.2222
$ws.Range("L79").Value = 3943.6
$ws.Range("M79").Value = -2440.2222
$ws.Range("N79").Value = -6127.6

# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 4168558
$ws.Range("I137").Value = 5883508.5
$ws.Range("J137").Value = 3678.5715
$ws.Range("K137").Value = 17650525.5
$ws.Range("L137").Value = 11035.7145
$ws.Range("M137").Value = -17647975.5
$ws.Range("N137").Value = -16135.7145

# Row 138: All-night Crafting (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 2827193
$ws.Range("I138").Value = 1456.125
$ws.Range("K138").Value = 4368.375
$ws.Range("M138").Value = 771.625


$ws = $wb.Worksheets.Item("ARM")

# Row 34: Insistent Sallets (Steel Sallet)
$ws.Range("H34").Value = 10001
$ws.Range("J34").Value = 10001
$ws.Range("L34").Value = 10001
$ws.Range("N34").Value = -10543

# Row 61: Dealing with the Tough Stuff (Cobalt Ingot)
$ws.Range("H61").Value = 43567356
$ws.Range("I61").Value = 100101620
$ws.Range("J61").Value = 79463
$ws.Range("K61").Value = 100101620
$ws.Range("L61").Value = 79463
$ws.Range("M61").Value = -100101408
$ws.Range("N61").Value = -79887

# Row 63: Rivets Run through It (Mythrite Rivets)
$ws.Range("H63").Value = 2953
$ws.Range("I63").Value = 2556.75
$ws.Range("J63").Value = 3481.3333
$ws.Range("K63").Value = 2556.75
$ws.Range("L63").Value = 3481.3333
$ws.Range("M63").Value = -1870.75
$ws.Range("N63").Value = -4853.3333

# Row 66: A Riveting Revival (L) (Mythrite Rivets)
$ws.Range("H66").Value = 2953
$ws.Range("I66").Value = 2556.75
$ws.Range("J66").Value = 3481.3333
$ws.Range("K66").Value = 12783.75
$ws.Range("L66").Value = 17406.6665
$ws.Range("M66").Value = -9351.75
$ws.Range("N66").Value = -24270.6665

# Row 132: Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws.Range("H132").Value = 107317.6
$ws.Range("I132").Value = 201902.4
$ws.Range("J132").Value = 75789.336
$ws.Range("K132").Value = 605707.2
$ws.Range("L132").Value = 227368.008
$ws.Range("M132").Value = -603177.2
$ws.Range("N132").Value = -232428.008

# Row 136: Metal with Mettle (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 43567356
$ws.Range("I136").Value = 100101620
$ws.Range("J136").Value = 79463
$ws.Range("K136").Value = 300304860
$ws.Range("L136").Value = 238389
$ws.Range("M136").Value = -300302310
$ws.Range("N136").Value = -243489


$ws = $wb.Worksheets.Item("BSM")

# Row 94: High Steal (High Steel Nugget)
$ws.Range("H94").Value = 697.4815
$ws.Range("I94").Value = 672.26086
$ws.Range("K94").Value = 672.26086
$ws.Range("M94").Value = -221.26086

# Row 105: Ingot to Wing It (Molybdenum Ingot)
$ws.Range("H105").Value = 21740846
$ws.Range("I105").Value = 25001666
$ws.Range("J105").Value = 2033.3334
$ws.Range("K105").Value = 25001666
$ws.Range("L105").Value = 2033.3334
$ws.Range("M105").Value = -24999919
$ws.Range("N105").Value = -5527.3334

# Row 107: The Gold Experience (Deepgold Nugget)
$ws.Range("H107").Value = 1544.8667
$ws.Range("I107").Value = 1629.4615
$ws.Range("J107").Value = 995
$ws.Range("K107").Value = 1629.4615
$ws.Range("L107").Value = 995
$ws.Range("M107").Value = 290.5385000000001
$ws.Range("N107").Value = -4835


$ws = $wb.Worksheets.Item("CRP")

# Row 105: Zelkova, My Love (Zelkova Lumber)
$ws.Range("H105").Value = 966.5
$ws.Range("I105").Value = 950
$ws.Range("K105").Value = 950
$ws.Range("M105").Value = 797

# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 48772.637
$ws.Range("I132").Value = 3388
$ws.Range("K132").Value = 10164
$ws.Range("M132").Value = -7634


$ws = $wb.Worksheets.Item("CUL")

# Row 113: Can't Eat Just One (Night Vinegar)
$ws.Range("H113").Value = 624.4138
$ws.Range("I113").Value = 538
$ws.Range("J113").Value = 657.3333
$ws.Range("K113").Value = 1614
$ws.Range("L113").Value = 1971.9999
$ws.Range("M113").Value = 556
$ws.Range("N113").Value = -6311.9999


$ws = $wb.Worksheets.Item("GSM")

# Row 70: Sky Is the Limit (Mythrite Ingot)
$ws.Range("H70").Value = 28448.906
$ws.Range("I70").Value = 43447.117
$ws.Range("K70").Value = 43447.117
$ws.Range("M70").Value = -43177.117

# Row 73: Hulls of Broken Dreams (L) (Mythrite Ingot)
$ws.Range("H73").Value = 28448.906
$ws.Range("I73").Value = 43447.117
$ws.Range("K73").Value = 43447.117
$ws.Range("M73").Value = -42511.117

# Row 80: Needs More Prayerbell (Hardsilver Ingot)
$ws.Range("H80").Value = 4020.1875
$ws.Range("I80").Value = 3561.6667
$ws.Range("J80").Value = 4126
$ws.Range("K80").Value = 3561.6667
$ws.Range("L80").Value = 4126
$ws.Range("M80").Value = -2563.6667
$ws.Range("N80").Value = -6122

# Row 83: With a Noise That Reaches Heaven (L) (Hardsilver Ingot)
$ws.Range("H83").Value = 4020.1875
$ws.Range("I83").Value = 3561.6667
$ws.Range("J83").Value = 4126
$ws.Range("K83").Value = 17808.3335
$ws.Range("L83").Value = 20630
$ws.Range("M83").Value = -12816.3335
$ws.Range("N83").Value = -30614


$ws = $wb.Worksheets.Item("LTW")

# Row 93: Hide to Go Seek (Gagana Leather)
$ws.Range("H93").Value = 1235.5264
$ws.Range("I93").Value = 1291.6154
$ws.Range("K93").Value = 1291.6154
$ws.Range("M93").Value = -43.61539999999991


$ws = $wb.Worksheets.Item("WVR")

# Row 96: Skills on Display (Ruby Cotton Cloth)
$ws.Range("H96").Value = 8812.25
$ws.Range("I96").Value = 3500
$ws.Range("J96").Value = 10583
$ws.Range("K96").Value = 3500
$ws.Range("L96").Value = 10583
$ws.Range("M96").Value = -2127
$ws.Range("N96").Value = -13329

# Row 122: Heavy Armoire (Dark Hempen Cloth)
$ws.Range("H122").Value = 3576.6667
$ws.Range("J122").Value = 4141.4287
$ws.Range("L122").Value = 12424.2861
$ws.Range("N122").Value = -17324.2861

# Row 132: Comfy Cabins (Snow Cotton Cloth)
$ws.Range("H132").Value = 48988.168
$ws.Range("I132").Value = 44554.176
$ws.Range("K132").Value = 133662.528
$ws.Range("M132").Value = -131132.528
